$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "01b93e67-7389-435b-9113-2936c4e1f3d3"
$ws.Range("B7").Value = "密毛魔芋"
$ws.Range("C7").Value = "有物室外"

$ws.Range("A8").Value = "ab514b5d-d30a-42e8-bf56-fe58063a5892"
$ws.Range("B8").Value = "綠背斜紋天蛾"
$ws.Range("C8").Value = "保全室"

$ws.Range("B7:C8").HorizontalAlignment = -4131
